$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing cell values (recalculated means) ---
$ws.Range("G2").Value = 1.81872300772161
$ws.Range("F9").Value = 0.01364
$ws.Range("G9").Value = 0.028919732795123
$ws.Range("F10").Value = 0.01364
$ws.Range("G10").Value = 0.028919732795123
$ws.Range("G11").Value = 0.59592708859653
$ws.Range("G12").Value = 0.59592708859653
$ws.Range("G13").Value = 0.5916923555310381
$ws.Range("G14").Value = 0.5916923555310381
$ws.Range("G19").Value = 1.69077438522579
$ws.Range("G26").Value = 0.0401788715679364
$ws.Range("G27").Value = 0.0401788715679364
$ws.Range("G28").Value = 0.5202542126491549
$ws.Range("G29").Value = 0.5202542126491549
$ws.Range("F30").Value = 0.5169
$ws.Range("G30").Value = 0.57004137850614
$ws.Range("F31").Value = 0.5169
$ws.Range("G31").Value = 0.57004137850614
$ws.Range("G36").Value = 1.60171778145221
$ws.Range("G43").Value = 0.0422350785190275
$ws.Range("G44").Value = 0.0422350785190275
$ws.Range("G45").Value = 0.5043785105409621
$ws.Range("G46").Value = 0.5043785105409621
$ws.Range("F47").Value = 0.51675
$ws.Range("G47").Value = 0.555773333333333
$ws.Range("N47").Value = 1.28778
$ws.Range("F48").Value = 0.51675
$ws.Range("G48").Value = 0.555773333333333
$ws.Range("N48").Value = 1.28778
$ws.Range("G60").Value = 0.0454981093860486
$ws.Range("G61").Value = 0.0454981093860486
$ws.Range("G62").Value = 0.5202764521497411
$ws.Range("G63").Value = 0.5202764521497411
$ws.Range("F64").Value = 0.5163
$ws.Range("G64").Value = 0.578448423430222
$ws.Range("L64").Value = 0.0563
$ws.Range("F65").Value = 0.5163
$ws.Range("G65").Value = 0.578448423430222
$ws.Range("L65").Value = 0.0563
$ws.Range("G78").Value = 0.0448782168152006
$ws.Range("G79").Value = 0.0448782168152006
$ws.Range("G80").Value = 0.528626452149741
$ws.Range("G81").Value = 0.528626452149741
$ws.Range("F82").Value = 0.5163
$ws.Range("G82").Value = 0.5896484234302219
$ws.Range("M82").Value = 1.12354
$ws.Range("F83").Value = 0.5163
$ws.Range("G83").Value = 0.5896484234302219
$ws.Range("M83").Value = 1.12354
$ws.Range("G96").Value = 0.0378636364860858
$ws.Range("G97").Value = 0.0378636364860858
$ws.Range("G98").Value = 0.536776452149741
$ws.Range("G99").Value = 0.536776452149741
$ws.Range("F100").Value = 0.46065
$ws.Range("G100").Value = 0.5956650900968889
$ws.Range("M100").Value = 1.19066
$ws.Range("F101").Value = 0.46065
$ws.Range("G101").Value = 0.5956650900968889
$ws.Range("M101").Value = 1.19066
$ws.Range("G115").Value = 0.0292467961869574
$ws.Range("L115").Value = 0.00389
$ws.Range("G116").Value = 0.0292467961869574
$ws.Range("L116").Value = 0.00389
$ws.Range("G117").Value = 0.558931452149741
$ws.Range("G118").Value = 0.558931452149741
$ws.Range("F120").Value = 0.5528999999999999
$ws.Range("G120").Value = 0.623093423430222
$ws.Range("M120").Value = 1.10639
$ws.Range("F121").Value = 0.5528999999999999
$ws.Range("G121").Value = 0.623093423430222
$ws.Range("M121").Value = 1.10639
$ws.Range("G135").Value = 0.0320555744385005
$ws.Range("L135").Value = 0.00555
$ws.Range("G136").Value = 0.0320555744385005
$ws.Range("L136").Value = 0.00555
$ws.Range("G137").Value = 0.559348118816407
$ws.Range("G138").Value = 0.559348118816407
$ws.Range("G140").Value = 0.627283423430222
$ws.Range("M140").Value = 1.09325
$ws.Range("G141").Value = 0.627283423430222
$ws.Range("M141").Value = 1.09325
$ws.Range("G155").Value = 0.0337217030503563
$ws.Range("G156").Value = 0.0337217030503563
$ws.Range("G175").Value = 0.03562762677938
$ws.Range("G176").Value = 0.03562762677938
$ws.Range("G195").Value = 0.0396247432038625
$ws.Range("G196").Value = 0.0396247432038625

# --- Append new data rows 206-224 ---
# Row 206
$ws.Range("A206").Value = "Mangaatua at d/s Woodville STP"
$ws.Range("B206").Value = "ASPM"
$ws.Range("C206").Value = "B"
$ws.Range("D206").Value = "2019 - 2023"
$ws.Range("E206").Value = "Impact"
$ws.Range("F206").Value = 0.418
$ws.Range("G206").Value = 0.4266
$ws.Range("H206").Value = 0.567
$ws.Range("I206").Value = 0.567
$ws.Range("J206").Value = ""
$ws.Range("K206").Value = ""
$ws.Range("L206").Value = 0.418
$ws.Range("M206").Value = 0.53235
$ws.Range("N206").Value = 0.567
$ws.Range("O206").Value = 1842194.8
$ws.Range("P206").Value = 5530097.413
$ws.Range("Q206").Value = "Tararua District"
$ws.Range("R206").Value = "Manawatū"
$ws.Range("S206").Value = "Upper Gorge"
$ws.Range("T206").Value = "Mana_9c"
$ws.Range("U206").Value = ""

# Row 207
$ws.Range("A207").Value = "Mangaatua at d/s Woodville STP"
$ws.Range("B207").Value = "DRP (95th Percentile)"
$ws.Range("C207").Value = "D"
$ws.Range("D207").Value = "2019 - 2023"
$ws.Range("E207").Value = "Impact"
$ws.Range("F207").Value = 0.03
$ws.Range("G207").Value = 0.0333559322033898
$ws.Range("H207").Value = 0.096
$ws.Range("I207").Value = 0.07235
$ws.Range("J207").Value = ""
$ws.Range("K207").Value = ""
$ws.Range("L207").Value = 0.035
$ws.Range("M207").Value = 0.04547
$ws.Range("N207").Value = 0.06012
$ws.Range("O207").Value = 1842194.8
$ws.Range("P207").Value = 5530097.413
$ws.Range("Q207").Value = "Tararua District"
$ws.Range("R207").Value = "Manawatū"
$ws.Range("S207").Value = "Upper Gorge"
$ws.Range("T207").Value = "Mana_9c"
$ws.Range("U207").Value = "mg/L"

# Row 208
$ws.Range("A208").Value = "Mangaatua at d/s Woodville STP"
$ws.Range("B208").Value = "DRP (Median)"
$ws.Range("C208").Value = "D"
$ws.Range("D208").Value = "2019 - 2023"
$ws.Range("E208").Value = "Impact"
$ws.Range("F208").Value = 0.03
$ws.Range("G208").Value = 0.0333559322033898
$ws.Range("H208").Value = 0.096
$ws.Range("I208").Value = 0.07235
$ws.Range("J208").Value = ""
$ws.Range("K208").Value = ""
$ws.Range("L208").Value = 0.035
$ws.Range("M208").Value = 0.04547
$ws.Range("N208").Value = 0.06012
$ws.Range("O208").Value = 1842194.8
$ws.Range("P208").Value = 5530097.413
$ws.Range("Q208").Value = "Tararua District"
$ws.Range("R208").Value = "Manawatū"
$ws.Range("S208").Value = "Upper Gorge"
$ws.Range("T208").Value = "Mana_9c"
$ws.Range("U208").Value = "mg/L"

# Row 209
$ws.Range("A209").Value = "Mangaatua at d/s Woodville STP"
$ws.Range("B209").Value = "E coli (>260)"
$ws.Range("C209").Value = "E"
$ws.Range("D209").Value = "2019 - 2023"
$ws.Range("E209").Value = "Impact"
$ws.Range("F209").Value = 420
$ws.Range("G209").Value = 1947.31091238693
$ws.Range("H209").Value = 16000
$ws.Range("I209").Value = 9700
$ws.Range("J209").Value = 44.0677966101695
$ws.Range("K209").Value = 72.8813559322034
$ws.Range("L209").Value = 592
$ws.Range("M209").Value = 3062.34
$ws.Range("N209").Value = 9073.88
$ws.Range("O209").Value = 1842194.8
$ws.Range("P209").Value = 5530097.413
$ws.Range("Q209").Value = "Tararua District"
$ws.Range("R209").Value = "Manawatū"
$ws.Range("S209").Value = "Upper Gorge"
$ws.Range("T209").Value = "Mana_9c"
$ws.Range("U209").Value = "% exceedances over 260/100 mL"

# Row 210
$ws.Range("A210").Value = "Mangaatua at d/s Woodville STP"
$ws.Range("B210").Value = "E coli (>540)"
$ws.Range("C210").Value = "E"
$ws.Range("D210").Value = "2019 - 2023"
$ws.Range("E210").Value = "Impact"
$ws.Range("F210").Value = 420
$ws.Range("G210").Value = 1947.31091238693
$ws.Range("H210").Value = 16000
$ws.Range("I210").Value = 9700
$ws.Range("J210").Value = 44.0677966101695
$ws.Range("K210").Value = 72.8813559322034
$ws.Range("L210").Value = 592
$ws.Range("M210").Value = 3062.34
$ws.Range("N210").Value = 9073.88
$ws.Range("O210").Value = 1842194.8
$ws.Range("P210").Value = 5530097.413
$ws.Range("Q210").Value = "Tararua District"
$ws.Range("R210").Value = "Manawatū"
$ws.Range("S210").Value = "Upper Gorge"
$ws.Range("T210").Value = "Mana_9c"
$ws.Range("U210").Value = "% exceedances over 540/100 mL"

# Row 211
$ws.Range("A211").Value = "Mangaatua at d/s Woodville STP"
$ws.Range("B211").Value = "E coli (Median)"
$ws.Range("C211").Value = "E"
$ws.Range("D211").Value = "2019 - 2023"
$ws.Range("E211").Value = "Impact"
$ws.Range("F211").Value = 420
$ws.Range("G211").Value = 1947.31091238693
$ws.Range("H211").Value = 16000
$ws.Range("I211").Value = 9700
$ws.Range("J211").Value = 44.0677966101695
$ws.Range("K211").Value = 72.8813559322034
$ws.Range("L211").Value = 592
$ws.Range("M211").Value = 3062.34
$ws.Range("N211").Value = 9073.88
$ws.Range("O211").Value = 1842194.8
$ws.Range("P211").Value = 5530097.413
$ws.Range("Q211").Value = "Tararua District"
$ws.Range("R211").Value = "Manawatū"
$ws.Range("S211").Value = "Upper Gorge"
$ws.Range("T211").Value = "Mana_9c"
$ws.Range("U211").Value = "E. coli/100 mL"

# Row 212
$ws.Range("A212").Value = "Mangaatua at d/s Woodville STP"
$ws.Range("B212").Value = "E coli (95th Percentile)"
$ws.Range("C212").Value = "E"
$ws.Range("D212").Value = "2019 - 2023"
$ws.Range("E212").Value = "Impact"
$ws.Range("F212").Value = 420
$ws.Range("G212").Value = 1947.31091238693
$ws.Range("H212").Value = 16000
$ws.Range("I212").Value = 9700
$ws.Range("J212").Value = 44.0677966101695
$ws.Range("K212").Value = 72.8813559322034
$ws.Range("L212").Value = 592
$ws.Range("M212").Value = 3062.34
$ws.Range("N212").Value = 9073.88
$ws.Range("O212").Value = 1842194.8
$ws.Range("P212").Value = 5530097.413
$ws.Range("Q212").Value = "Tararua District"
$ws.Range("R212").Value = "Manawatū"
$ws.Range("S212").Value = "Upper Gorge"
$ws.Range("T212").Value = "Mana_9c"
$ws.Range("U212").Value = "E. coli/100 mL"

# Row 213
$ws.Range("A213").Value = "Mangaatua at d/s Woodville STP"
$ws.Range("B213").Value = "MCI"
$ws.Range("C213").Value = "C"
$ws.Range("D213").Value = "2019 - 2023"
$ws.Range("E213").Value = "Impact"
$ws.Range("F213").Value = 107.27
$ws.Range("G213").Value = 109.988
$ws.Range("H213").Value = 126.67
$ws.Range("I213").Value = 126.67
$ws.Range("J213").Value = ""
$ws.Range("K213").Value = ""
$ws.Range("L213").Value = 107.27
$ws.Range("M213").Value = 124.6855
$ws.Range("N213").Value = 126.67
$ws.Range("O213").Value = 1842194.8
$ws.Range("P213").Value = 5530097.413
$ws.Range("Q213").Value = "Tararua District"
$ws.Range("R213").Value = "Manawatū"
$ws.Range("S213").Value = "Upper Gorge"
$ws.Range("T213").Value = "Mana_9c"
$ws.Range("U213").Value = ""

# Row 214
$ws.Range("A214").Value = "Mangaatua at d/s Woodville STP"
$ws.Range("B214").Value = "Ammoniacal-N (95th Percentile)"
$ws.Range("C214").Value = "B"
$ws.Range("D214").Value = "2019 - 2023"
$ws.Range("E214").Value = "Impact"
$ws.Range("F214").Value = 0.0345
$ws.Range("G214").Value = 0.048521321816706
$ws.Range("H214").Value = 0.309482048166445
$ws.Range("I214").Value = 0.11325
$ws.Range("J214").Value = ""
$ws.Range("K214").Value = ""
$ws.Range("L214").Value = 0.02749
$ws.Range("M214").Value = 0.07431
$ws.Range("N214").Value = 0.09335
$ws.Range("O214").Value = 1842194.8
$ws.Range("P214").Value = 5530097.413
$ws.Range("Q214").Value = "Tararua District"
$ws.Range("R214").Value = "Manawatū"
$ws.Range("S214").Value = "Upper Gorge"
$ws.Range("T214").Value = "Mana_9c"
$ws.Range("U214").Value = "mg NH4-N/L"

# Row 215
$ws.Range("A215").Value = "Mangaatua at d/s Woodville STP"
$ws.Range("B215").Value = "Ammoniacal-N (Median)"
$ws.Range("C215").Value = "B"
$ws.Range("D215").Value = "2019 - 2023"
$ws.Range("E215").Value = "Impact"
$ws.Range("F215").Value = 0.0345
$ws.Range("G215").Value = 0.048521321816706
$ws.Range("H215").Value = 0.309482048166445
$ws.Range("I215").Value = 0.11325
$ws.Range("J215").Value = ""
$ws.Range("K215").Value = ""
$ws.Range("L215").Value = 0.02749
$ws.Range("M215").Value = 0.07431
$ws.Range("N215").Value = 0.09335
$ws.Range("O215").Value = 1842194.8
$ws.Range("P215").Value = 5530097.413
$ws.Range("Q215").Value = "Tararua District"
$ws.Range("R215").Value = "Manawatū"
$ws.Range("S215").Value = "Upper Gorge"
$ws.Range("T215").Value = "Mana_9c"
$ws.Range("U215").Value = "mg NH4-N/L"

# Row 216
$ws.Range("A216").Value = "Mangaatua at d/s Woodville STP"
$ws.Range("B216").Value = "Nitrate-N (95th Percentile)"
$ws.Range("C216").Value = "B"
$ws.Range("D216").Value = "2019 - 2023"
$ws.Range("E216").Value = "Impact"
$ws.Range("F216").Value = 0.348
$ws.Range("G216").Value = 0.480829273114289
$ws.Range("H216").Value = 1.86
$ws.Range("I216").Value = 1.5075
$ws.Range("J216").Value = ""
$ws.Range("K216").Value = ""
$ws.Range("L216").Value = 0.122
$ws.Range("M216").Value = 0.79441
$ws.Range("N216").Value = 1.0278
$ws.Range("O216").Value = 1842194.8
$ws.Range("P216").Value = 5530097.413
$ws.Range("Q216").Value = "Tararua District"
$ws.Range("R216").Value = "Manawatū"
$ws.Range("S216").Value = "Upper Gorge"
$ws.Range("T216").Value = "Mana_9c"
$ws.Range("U216").Value = "mg NO3-N/L"

# Row 217
$ws.Range("A217").Value = "Mangaatua at d/s Woodville STP"
$ws.Range("B217").Value = "Nitrate-N (Median)"
$ws.Range("C217").Value = "A"
$ws.Range("D217").Value = "2019 - 2023"
$ws.Range("E217").Value = "Impact"
$ws.Range("F217").Value = 0.348
$ws.Range("G217").Value = 0.480829273114289
$ws.Range("H217").Value = 1.86
$ws.Range("I217").Value = 1.5075
$ws.Range("J217").Value = ""
$ws.Range("K217").Value = ""
$ws.Range("L217").Value = 0.122
$ws.Range("M217").Value = 0.79441
$ws.Range("N217").Value = 1.0278
$ws.Range("O217").Value = 1842194.8
$ws.Range("P217").Value = 5530097.413
$ws.Range("Q217").Value = "Tararua District"
$ws.Range("R217").Value = "Manawatū"
$ws.Range("S217").Value = "Upper Gorge"
$ws.Range("T217").Value = "Mana_9c"
$ws.Range("U217").Value = "mg NO3-N/L"

# Row 218
$ws.Range("A218").Value = "Mangaatua at d/s Woodville STP"
$ws.Range("B218").Value = "QMCI"
$ws.Range("C218").Value = "B"
$ws.Range("D218").Value = "2019 - 2023"
$ws.Range("E218").Value = "Impact"
$ws.Range("F218").Value = 6.06
$ws.Range("G218").Value = 6.1312
$ws.Range("H218").Value = 7.2
$ws.Range("I218").Value = 7.2
$ws.Range("J218").Value = ""
$ws.Range("K218").Value = ""
$ws.Range("L218").Value = 6.06
$ws.Range("M218").Value = 6.9781
$ws.Range("N218").Value = 7.2
$ws.Range("O218").Value = 1842194.8
$ws.Range("P218").Value = 5530097.413
$ws.Range("Q218").Value = "Tararua District"
$ws.Range("R218").Value = "Manawatū"
$ws.Range("S218").Value = "Upper Gorge"
$ws.Range("T218").Value = "Mana_9c"
$ws.Range("U218").Value = ""

# Row 219
$ws.Range("A219").Value = "Mangaatua at d/s Woodville STP"
$ws.Range("B219").Value = "Soluble Inorganic Nitrogen (95th Percentile)"
$ws.Range("C219").Value = ""
$ws.Range("D219").Value = "2019 - 2023"
$ws.Range("E219").Value = "Impact"
$ws.Range("F219").Value = 0.492
$ws.Range("G219").Value = 0.579305084745763
$ws.Range("H219").Value = 1.956
$ws.Range("I219").Value = 1.6278
$ws.Range("J219").Value = ""
$ws.Range("K219").Value = ""
$ws.Range("L219").Value = 0.201
$ws.Range("M219").Value = 0.95205
$ws.Range("N219").Value = 1.12406
$ws.Range("O219").Value = 1842194.8
$ws.Range("P219").Value = 5530097.413
$ws.Range("Q219").Value = "Tararua District"
$ws.Range("R219").Value = "Manawatū"
$ws.Range("S219").Value = "Upper Gorge"
$ws.Range("T219").Value = "Mana_9c"
$ws.Range("U219").Value = "g/m3"

# Row 220
$ws.Range("A220").Value = "Mangaatua at d/s Woodville STP"
$ws.Range("B220").Value = "Soluble Inorganic Nitrogen (Median)"
$ws.Range("C220").Value = ""
$ws.Range("D220").Value = "2019 - 2023"
$ws.Range("E220").Value = "Impact"
$ws.Range("F220").Value = 0.492
$ws.Range("G220").Value = 0.579305084745763
$ws.Range("H220").Value = 1.956
$ws.Range("I220").Value = 1.6278
$ws.Range("J220").Value = ""
$ws.Range("K220").Value = ""
$ws.Range("L220").Value = 0.201
$ws.Range("M220").Value = 0.95205
$ws.Range("N220").Value = 1.12406
$ws.Range("O220").Value = 1842194.8
$ws.Range("P220").Value = 5530097.413
$ws.Range("Q220").Value = "Tararua District"
$ws.Range("R220").Value = "Manawatū"
$ws.Range("S220").Value = "Upper Gorge"
$ws.Range("T220").Value = "Mana_9c"
$ws.Range("U220").Value = "g/m3"

# Row 221
$ws.Range("A221").Value = "Mangaatua at d/s Woodville STP"
$ws.Range("B221").Value = "Total Nitrogen (95th Percentile)"
$ws.Range("C221").Value = ""
$ws.Range("D221").Value = "2019 - 2023"
$ws.Range("E221").Value = "Impact"
$ws.Range("F221").Value = 0.89
$ws.Range("G221").Value = 0.952033898305085
$ws.Range("H221").Value = 2.72
$ws.Range("I221").Value = 2.094
$ws.Range("J221").Value = ""
$ws.Range("K221").Value = ""
$ws.Range("L221").Value = 0.465
$ws.Range("M221").Value = 1.4105
$ws.Range("N221").Value = 1.8656
$ws.Range("O221").Value = 1842194.8
$ws.Range("P221").Value = 5530097.413
$ws.Range("Q221").Value = "Tararua District"
$ws.Range("R221").Value = "Manawatū"
$ws.Range("S221").Value = "Upper Gorge"
$ws.Range("T221").Value = "Mana_9c"
$ws.Range("U221").Value = "g/m3"

# Row 222
$ws.Range("A222").Value = "Mangaatua at d/s Woodville STP"
$ws.Range("B222").Value = "Total Nitrogen (Median)"
$ws.Range("C222").Value = ""
$ws.Range("D222").Value = "2019 - 2023"
$ws.Range("E222").Value = "Impact"
$ws.Range("F222").Value = 0.89
$ws.Range("G222").Value = 0.952033898305085
$ws.Range("H222").Value = 2.72
$ws.Range("I222").Value = 2.094
$ws.Range("J222").Value = ""
$ws.Range("K222").Value = ""
$ws.Range("L222").Value = 0.465
$ws.Range("M222").Value = 1.4105
$ws.Range("N222").Value = 1.8656
$ws.Range("O222").Value = 1842194.8
$ws.Range("P222").Value = 5530097.413
$ws.Range("Q222").Value = "Tararua District"
$ws.Range("R222").Value = "Manawatū"
$ws.Range("S222").Value = "Upper Gorge"
$ws.Range("T222").Value = "Mana_9c"
$ws.Range("U222").Value = "g/m3"

# Row 223
$ws.Range("A223").Value = "Mangaatua at d/s Woodville STP"
$ws.Range("B223").Value = "Total Phosphorus (95th Percentile)"
$ws.Range("C223").Value = ""
$ws.Range("D223").Value = "2019 - 2023"
$ws.Range("E223").Value = "Impact"
$ws.Range("F223").Value = 0.068
$ws.Range("G223").Value = 0.102220338983051
$ws.Range("H223").Value = 0.539
$ws.Range("I223").Value = 0.3382
$ws.Range("J223").Value = ""
$ws.Range("K223").Value = ""
$ws.Range("L223").Value = 0.0765
$ws.Range("M223").Value = 0.15546
$ws.Range("N223").Value = 0.17658
$ws.Range("O223").Value = 1842194.8
$ws.Range("P223").Value = 5530097.413
$ws.Range("Q223").Value = "Tararua District"
$ws.Range("R223").Value = "Manawatū"
$ws.Range("S223").Value = "Upper Gorge"
$ws.Range("T223").Value = "Mana_9c"
$ws.Range("U223").Value = "g/m3"

# Row 224
$ws.Range("A224").Value = "Mangaatua at d/s Woodville STP"
$ws.Range("B224").Value = "Total Phosphorus (Median)"
$ws.Range("C224").Value = ""
$ws.Range("D224").Value = "2019 - 2023"
$ws.Range("E224").Value = "Impact"
$ws.Range("F224").Value = 0.068
$ws.Range("G224").Value = 0.102220338983051
$ws.Range("H224").Value = 0.539
$ws.Range("I224").Value = 0.3382
$ws.Range("J224").Value = ""
$ws.Range("K224").Value = ""
$ws.Range("L224").Value = 0.0765
$ws.Range("M224").Value = 0.15546
$ws.Range("N224").Value = 0.17658
$ws.Range("O224").Value = 1842194.8
$ws.Range("P224").Value = 5530097.413
$ws.Range("Q224").Value = "Tararua District"
$ws.Range("R224").Value = "Manawatū"
$ws.Range("S224").Value = "Upper Gorge"
$ws.Range("T224").Value = "Mana_9c"
$ws.Range("U224").Value = "g/m3"

